$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()
